$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")
$ws.Activate()

# Update ventilation properties (column H, rows 26:181) from 0.9 to 0.5
$ws.Range("H26:H181").Value = 0.5

# Update the active selection to match the edited range
$ws.Range("H26:H181").Select()
